$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 13.26539866666667
$ws.Cells.Item(2, 8).Value = 39.796196
$ws.Cells.Item(2, 9).Value = 0.164744212542501
$ws.Cells.Item(2, 10).Value = 0.164744212542501
$ws.Cells.Item(2, 13).Value = 0.6340813333333334
$ws.Cells.Item(2, 14).Value = 1.902244
$ws.Cells.Item(2, 15).Value = 0.003487630722257058
$ws.Cells.Item(2, 16).Value = 0.003487630722257058
$ws.Cells.Item(2, 17).Value = 8.411341673758225
$ws.Cells.Item(2, 18).Value = 75.70207506382401
$ws.Cells.Item(2, 19).Value = 0.000574566976977273
$ws.Cells.Item(2, 20).Value = 0.0005745669769772729

$ws.Cells.Item(3, 7).Value = 13.26539866666667
$ws.Cells.Item(3, 8).Value = 39.796196
$ws.Cells.Item(3, 9).Value = 0.164744212542501
$ws.Cells.Item(3, 10).Value = 0.164744212542501
$ws.Cells.Item(3, 15).Value = 0.8644503444376447
$ws.Cells.Item(3, 16).Value = 0.8644503444376448
$ws.Cells.Item(3, 17).Value = 2084.8500848041
$ws.Cells.Item(3, 18).Value = 18763.6507632369
$ws.Cells.Item(3, 19).Value = 0.1424131912764735
$ws.Cells.Item(3, 20).Value = 0.1424131912764735

$ws.Cells.Item(4, 7).Value = 13.26539866666667
$ws.Cells.Item(4, 8).Value = 39.796196
$ws.Cells.Item(4, 9).Value = 0.164744212542501
$ws.Cells.Item(4, 10).Value = 0.164744212542501
$ws.Cells.Item(4, 13).Value = 24.01001466666667
$ws.Cells.Item(4, 14).Value = 72.030044
$ws.Cells.Item(4, 15).Value = 0.1320620248400982
$ws.Cells.Item(4, 16).Value = 0.1320620248400982
$ws.Cells.Item(4, 17).Value = 318.5024165458472
$ws.Cells.Item(4, 18).Value = 2866.521748912624
$ws.Cells.Item(4, 19).Value = 0.02175645428905018
$ws.Cells.Item(4, 20).Value = 0.02175645428905018

$ws.Cells.Item(5, 9).Value = 0.6678031736949381
$ws.Cells.Item(5, 10).Value = 0.6678031736949381
$ws.Cells.Item(5, 13).Value = 0.6340813333333334
$ws.Cells.Item(5, 14).Value = 1.902244
$ws.Cells.Item(5, 15).Value = 0.003487630722257058
$ws.Cells.Item(5, 16).Value = 0.003487630722257058
$ws.Cells.Item(5, 17).Value = 34.09601210312089
$ws.Cells.Item(5, 18).Value = 306.864108928088
$ws.Cells.Item(5, 19).Value = 0.002329050864999233
$ws.Cells.Item(5, 20).Value = 0.002329050864999232

$ws.Cells.Item(6, 9).Value = 0.6678031736949381
$ws.Cells.Item(6, 10).Value = 0.6678031736949381
$ws.Cells.Item(6, 15).Value = 0.8644503444376447
$ws.Cells.Item(6, 16).Value = 0.8644503444376448
$ws.Cells.Item(6, 18).Value = 76059.88299322156
$ws.Cells.Item(6, 19).Value = 0.5772826835171415
$ws.Cells.Item(6, 20).Value = 0.5772826835171416

$ws.Cells.Item(7, 9).Value = 0.6678031736949381
$ws.Cells.Item(7, 10).Value = 0.6678031736949381
$ws.Cells.Item(7, 13).Value = 24.01001466666667
$ws.Cells.Item(7, 14).Value = 72.030044
$ws.Cells.Item(7, 15).Value = 0.1320620248400982
$ws.Cells.Item(7, 16).Value = 0.1320620248400982
$ws.Cells.Item(7, 17).Value = 1291.073727667076
$ws.Cells.Item(7, 18).Value = 11619.66354900369
$ws.Cells.Item(7, 19).Value = 0.0881914393127973
$ws.Cells.Item(7, 20).Value = 0.0881914393127973

$ws.Cells.Item(8, 7).Value = 13.48348233333333
$ws.Cells.Item(8, 8).Value = 40.450447
$ws.Cells.Item(8, 9).Value = 0.1674526137625609
$ws.Cells.Item(8, 10).Value = 0.1674526137625609
$ws.Cells.Item(8, 13).Value = 0.6340813333333334
$ws.Cells.Item(8, 14).Value = 1.902244
$ws.Cells.Item(8, 15).Value = 0.003487630722257058
$ws.Cells.Item(8, 16).Value = 0.003487630722257058
$ws.Cells.Item(8, 17).Value = 8.549624455896446
$ws.Cells.Item(8, 18).Value = 76.946620103068
$ws.Cells.Item(8, 19).Value = 0.0005840128802805525
$ws.Cells.Item(8, 20).Value = 0.0005840128802805523

$ws.Cells.Item(9, 7).Value = 13.48348233333333
$ws.Cells.Item(9, 8).Value = 40.450447
$ws.Cells.Item(9, 9).Value = 0.1674526137625609
$ws.Cells.Item(9, 10).Value = 0.1674526137625609
$ws.Cells.Item(9, 15).Value = 0.8644503444376447
$ws.Cells.Item(9, 16).Value = 0.8644503444376448
$ws.Cells.Item(9, 17).Value = 2119.125100758719
$ws.Cells.Item(9, 18).Value = 19072.12590682847
$ws.Cells.Item(9, 19).Value = 0.1447544696440297
$ws.Cells.Item(9, 20).Value = 0.1447544696440297

$ws.Cells.Item(10, 7).Value = 13.48348233333333
$ws.Cells.Item(10, 8).Value = 40.450447
$ws.Cells.Item(10, 9).Value = 0.1674526137625609
$ws.Cells.Item(10, 10).Value = 0.1674526137625609
$ws.Cells.Item(10, 13).Value = 24.01001466666667
$ws.Cells.Item(10, 14).Value = 72.030044
$ws.Cells.Item(10, 15).Value = 0.1320620248400982
$ws.Cells.Item(10, 16).Value = 0.1320620248400982
$ws.Cells.Item(10, 17).Value = 323.7386085810742
$ws.Cells.Item(10, 18).Value = 2913.647477229668
$ws.Cells.Item(10, 19).Value = 0.02211413123825068
$ws.Cells.Item(10, 20).Value = 0.02211413123825068
